{"js": "const body = context.document.body;\nconst results = body.search(\"RPC Explorer\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Insight Explorer\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"RPC Explorer\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Insight Explorer\"\n$find.Execute(\n    [ref]\"RPC Explorer\",   # FindText\n    [ref]$true,             # MatchCase\n    [ref]$false,            # MatchWholeWord\n    [ref]$false,            # MatchWildcards\n    [ref]$false,            # MatchSoundsLike\n    [ref]$false,            # MatchAllWordForms\n    [ref]$true,             # Forward\n    1,                      # Wrap = wdFindContinue\n    [ref]$false,            # Format\n    [ref]\"Insight Explorer\",# ReplaceWith\n    2                       # Replace = wdReplaceAll\n)\n"}
